$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 7 (the "d=7" row), shifting the
# existing "d=7" and "d=10" rows down to rows 8 and 9.
$ws.Rows("7:7").Insert()

# Copy the style of the label cell in the row above (A6, which holds "d=5")
# onto the new label cell A7 so it keeps the same formatting (bold, border,
# centered) as the rest of column A.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new row with the "d=6" label and its corresponding values.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.95031541226649
$ws.Range("C7").Value = 98.03180981637902
$ws.Range("D7").Value = 98.07944176787048
$ws.Range("E7").Value = 98.0267340907819

$wb.Save()
